$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing B3/C3 values down into column A of rows 4 and 5
$ws.Range("A4").Value = 2.2
$ws.Range("A5").Value = 3.3

# Fill in column B values for rows 3-5
$ws.Range("B3").Value = 11.1
$ws.Range("B4").Value = 12.2
$ws.Range("B5").Value = 13.3

# Fill in column C values for rows 3-5
$ws.Range("C3").Value = 21.1
$ws.Range("C4").Value = 22.2
$ws.Range("C5").Value = 23.3
